# [Kadastro App] Yeni kayit eklendi: 3020
# Appends the new record row (row 79) to both the "Kayitlar" master sheet
# and the "Erdemli" per-unit sheet, mirroring the existing rows exactly
# (all columns stored as literal text, including the numeric-looking and
# date-looking ones).

$wb = $excel.ActiveWorkbook

$newRow = @("3020", "2025-09-11", "Erdemli", "1", "3B", "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowIndex = $ws.UsedRange.Rows.Count + 1

    $target = $ws.Range($ws.Cells.Item($rowIndex, 1), $ws.Cells.Item($rowIndex, 6))

    # Force every cell in the new row to be treated as plain text so Excel
    # doesn't auto-coerce "3020" to a number or "2025-09-11" to a date --
    # matching how the rest of the sheet stores these values as strings.
    $target.NumberFormat = "@"

    for ($i = 0; $i -lt $newRow.Length; $i++) {
        $ws.Cells.Item($rowIndex, $i + 1).Value = $newRow[$i]
    }

    # Drop the explicit "text" number-format override again so the new
    # cells end up with the same (default) style as their neighbours.
    $target.ClearFormats()
}
